$d = $word.ActiveDocument

# Target paragraph: "This is a Microsoft word document."
$para = $d.Paragraphs(1)
$pRange = $para.Range

# Pull the paragraph's real opening <w:p ...> tag (paraId/rsid/etc.) so the
# rebuilt paragraph keeps its identity instead of minting a new one.
$wx = $pRange.WordOpenXML
if ($wx -match "<w:p[ >][^>]*>") {
    $pOpenTag = $matches[0]
} else {
    $pOpenTag = "<w:p>"
}

# Existing run text for this paragraph (kept byte-for-byte as the first run).
$existingText = "This is a Microsoft word document."

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
if ($pOpenTag -notmatch 'xmlns:w14') {
    $wNs = $wNs + ' xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
}

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document ' + $wNs + '>' +
       '<w:body>' + $pOpenTag +
       '<w:r><w:t>' + $existingText + '</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
       '<w:r><w:t>Changed main</w:t></w:r>' +
       '<w:r><w:t>)</w:t></w:r>' +
       '</w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

[void]$pRange.InsertXML($xml)
